# Apply updated cryptocurrency price/volume/coin data to Sheet1.
# Column D (Price) values are prefixed with a leading apostrophe so Excel
# keeps them as literal text (e.g. "60.144.14", "1.00") instead of coercing
# them into numbers, matching the existing inline-string text cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''60.144.14'
$ws.Range('E2').Value = '  -2.98%  '
# Row 3
$ws.Range('D3').Value = '''3.298.02'
$ws.Range('E3').Value = '  -3.61%  '
# Row 4
$ws.Range('E4').Value = '  +0.04%  '
# Row 5
$ws.Range('D5').Value = '''556.93'
$ws.Range('E5').Value = '  -3.84%  '
# Row 6
$ws.Range('D6').Value = '''141.23'
$ws.Range('E6').Value = '  -8.22%  '
# Row 7
$ws.Range('E7').Value = '  +0.03%  '
# Row 8
$ws.Range('D8').Value = '''3.298.30'
$ws.Range('E8').Value = '  -3.61%  '
# Row 9
$ws.Range('E9').Value = '  -3.75%  '
# Row 10
$ws.Range('D10').Value = '''7.89'
$ws.Range('E10').Value = '  -1.68%  '
# Row 11
$ws.Range('E11').Value = '  -5.29%  '
# Row 12
$ws.Range('D12').Value = '''0.408'
$ws.Range('E12').Value = '  -2.59%  '
# Row 13
$ws.Range('D13').Value = '''3.862.03'
$ws.Range('E13').Value = '  -3.71%  '
# Row 14
$ws.Range('E14').Value = '  +0.02%  '
# Row 15
$ws.Range('D15').Value = '''26.76'
$ws.Range('E15').Value = '  -5.61%  '
# Row 16
$ws.Range('D16').Value = '''3.289.74'
$ws.Range('E16').Value = '  -3.73%  '
# Row 17
$ws.Range('E17').Value = '  -5.12%  '
# Row 18
$ws.Range('D18').Value = '''60.199.81'
$ws.Range('E18').Value = '  -2.93%  '
# Row 19
$ws.Range('D19').Value = '''6.06'
$ws.Range('E19').Value = '  -8.03%  '
# Row 20
$ws.Range('D20').Value = '''13.74'
$ws.Range('E20').Value = '  -5.13%  '
# Row 21
$ws.Range('D21').Value = '''8.55'
$ws.Range('E21').Value = '  -4.44%  '
# Row 22
$ws.Range('D22').Value = '''372.96'
$ws.Range('E22').Value = '  -2.57%  '
# Row 23
$ws.Range('E23').Value = '  +0.12%  '
# Row 24
$ws.Range('D24').Value = '''72.58'
$ws.Range('E24').Value = '  -4.23%  '
# Row 25
$ws.Range('E25').Value = '  -6.92%  '
# Row 26
$ws.Range('D26').Value = '''3.424.86'
$ws.Range('E26').Value = '  -3.82%  '
# Row 27
$ws.Range('D27').Value = '''0.0000102'
$ws.Range('E27').Value = '  -9.62%  '
# Row 28
$ws.Range('D28').Value = '''0.173'
$ws.Range('E28').Value = '  -3.40%  '
# Row 29
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.46%  '
# Row 30
$ws.Range('D30').Value = '''7.05'
$ws.Range('E30').Value = '  -7.86%  '
# Row 31
$ws.Range('E31').Value = '  -0.07%  '
# Row 32
$ws.Range('D32').Value = '''2.02'
$ws.Range('E32').Value = '  -4.97%  '
# Row 33
$ws.Range('D33').Value = '''7.44'
$ws.Range('E33').Value = '  -5.55%  '
# Row 34
$ws.Range('D34').Value = '''22.55'
$ws.Range('E34').Value = '  -3.16%  '
# Row 35
$ws.Range('E35').Value = '  -7.70%  '
# Row 36
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = '''166.48'
$ws.Range('E36').Value = '  -0.96%  '
# Row 37
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = '''5.05'
$ws.Range('E37').Value = '  -9.17%  '
# Row 38
$ws.Range('D38').Value = '''1.53'
$ws.Range('E38').Value = '  -4.88%  '
# Row 39
$ws.Range('D39').Value = '''6.62'
$ws.Range('E39').Value = '  -5.10%  '
# Row 40
$ws.Range('D40').Value = '''3.329.40'
$ws.Range('E40').Value = '  -3.69%  '
# Row 41
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').Value = '''0.0722'
$ws.Range('E41').Value = '  -8.14%  '
# Row 42
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '''25.63'
$ws.Range('E42').Value = '  -17.67%  '
# Row 43
$ws.Range('D43').Value = '''41.59'
$ws.Range('E43').Value = '  -2.60%  '
# Row 44
$ws.Range('D44').Value = '''0.748'
$ws.Range('E44').Value = '  -4.25%  '
# Row 45
$ws.Range('D45').Value = '''1.11'
$ws.Range('E45').Value = '  -4.28%  '
# Row 46
$ws.Range('D46').Value = '''4.09'
$ws.Range('E46').Value = '  -7.75%  '
# Row 47
$ws.Range('E47').Value = '  -7.03%  '
# Row 48
$ws.Range('E48').Value = '  -0.04%  '
# Row 49
$ws.Range('D49').Value = '''2.316.25'
$ws.Range('E49').Value = '  -9.24%  '
# Row 50
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '''21.56'
$ws.Range('E50').Value = '  -6.40%  '
# Row 51
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').Value = '''6.35'
$ws.Range('E51').Value = '  -6.97%  '
